# Update the "Förändrad" (Changed) date column (column C) for all data rows.
# Every value in C2:C472 moves forward by one day (date serial 45189 -> 45190).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = $ws.UsedRange.Rows.Count
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value2
    if ($current -ne $null -and $current -ne "") {
        $cell.Value2 = $current + 1
    }
}
